$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the rows for years 2000, 2005, 2006, 2007, 2008, 2009 (original rows 2-7),
# shifting the cells below upward (xlShiftUp = -4162).
# This moves all subsequent rows (originally 8-18, years 2010-2020) up to rows 2-12.
$ws.Range("A2:D7").Delete(-4162)

# Append the two new years (2021 and 2022) after the shifted data, which now ends at row 12.
$ws.Range("A13").Value = "2021年"
$ws.Range("B13").Value = 1152762.31
$ws.Range("C13").Value = 2170486.01
$ws.Range("D13").Value = 1017723.7

$ws.Range("A14").Value = "2022年"
$ws.Range("B14").Value = 1123330.223849
$ws.Range("C14").Value = 2076367.760936
$ws.Range("D14").Value = 953037.537087

# Match the formatting used by the existing year-label cells (column A) by copying
# just the formatting (xlPasteFormats = -4122) from the row above onto the new rows.
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null
